$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement section
$ws.Range("F10").Value = 2067100
$ws.Range("F15").Value = 122900
$ws.Range("F17").Value = 4954800
$ws.Range("F18").Value = 235800
$ws.Range("F20").Value = -5300
$ws.Range("J21").Value = "NA"
$ws.Range("F24").Value = 84000
$ws.Range("F26").Value = 142900
$ws.Range("F27").Value = 142900
$ws.Range("F32").Value = 5300
$ws.Range("F33").Value = 142700
$ws.Range("F35").Value = 142700

# Balance Sheet section
$ws.Range("D58").Value = 4800
$ws.Range("E58").Value = 5500
$ws.Range("F58").Value = 6000
$ws.Range("G58").Value = 4000
$ws.Range("H58").Value = 900
$ws.Range("I58").Value = 1300
$ws.Range("J58").Value = "NA"

$ws.Range("D59").Value = 233600
$ws.Range("E59").Value = 272600
$ws.Range("F59").Value = 290200
$ws.Range("G59").Value = 287800
$ws.Range("H59").Value = 210700
$ws.Range("I59").Value = 234700

$ws.Range("D61").Value = 211000
$ws.Range("E61").Value = 121200
$ws.Range("F61").Value = 81700
$ws.Range("H61").Value = 77100
$ws.Range("I61").Value = 172000

$ws.Range("D62").Value = 181500
$ws.Range("E62").Value = 157300
$ws.Range("F62").Value = 159700
$ws.Range("G62").Value = 237900
$ws.Range("H62").Value = 183600
$ws.Range("I62").Value = 193800

# Cash Flow Statement section
$ws.Range("F81").Value = 142700
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"
